$d = $word.ActiveDocument

# The target paragraph is the one that still talks about the constellation
# "Perseus" (spread across four runs: "...Sternbildes ", "Perseus", " ",
# "am Nachthimmel..."). We replace the whole paragraph's text with a single,
# unformatted run that names "Zwillinge" instead, merging the four runs into
# one run with no explicit run formatting (matching the target edit).

$newText = "Mach mit an einer weltweiten Kampagne, die schwächsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Zwillinge am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Sternbildes*Perseus*Nachthimmel*") {
        # Range covering the paragraph's text, excluding the trailing pilcrow.
        $target = $d.Range($p.Range.Start, $p.Range.End - 1)

        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
               '</w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'

        $target.InsertXML($xml)
    }
}
